# Upload stage 3 evidences B1-B2
$wb = $excel.ActiveWorkbook

# --- Fill in B1 evidence (A2/A3) ---
$wsB1 = $wb.Worksheets.Item("B1")
$wsB1.Range("A2").Value = "FBFD588EDFD9275A1FFF4F141BA4827F83559F548168341C39EDAB0BB8E1241A"
$wsB1.Range("A3").Value = "E331FC592E2A72E0D60590C5B767CE737AEB89E097AEDDF1FEBAF9DB2B2069DC"
$wsB1.Activate()
$wsB1.Range("A4").Select()

# --- Fill in B2 evidence (A2/A3) ---
$wsB2 = $wb.Worksheets.Item("B2")
$wsB2.Range("A2").Value = "3161AEBB129ECB3D4A4E747DC5F1DA2EA4AA0FA3C7A4734D6454F82ABCBDB94A"
$wsB2.Range("A3").Value = "8275A3624278481BA0F4161865AD61EBF95FB656CDF6B7E7D0BDB35C6D841237"
$wsB2.Activate()
$wsB2.Range("A4").Select()
